$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.93002969492931
$ws.Range("C2").Value = 8.409604367810166
$ws.Range("E2").Value = 14.1873705680483
$ws.Range("F2").Value = 44.21552489145481
$ws.Range("G2").Value = 48.60556165666967
$ws.Range("H2").Value = 19.3085868293916
$ws.Range("J2").Value = 9.409833416744263
$ws.Range("K2").Value = 11.79130365214296
$ws.Range("L2").Value = 11.61897990344499
$ws.Range("M2").Value = 16.79059600469641
$ws.Range("N2").Value = 22.74105422421475
$ws.Range("B3").Value = 15.7817465609243
$ws.Range("C3").Value = 8.377977149994367
$ws.Range("E3").Value = 14.20153077677018
$ws.Range("F3").Value = 44.23310394848778
$ws.Range("G3").Value = 48.65555199953246
$ws.Range("H3").Value = 19.35373911895653
$ws.Range("J3").Value = 9.413562132509229
$ws.Range("K3").Value = 11.68717809337669
$ws.Range("L3").Value = 11.62263412457405
$ws.Range("M3").Value = 16.77372773054288
$ws.Range("N3").Value = 22.80389720399781
$ws.Range("B4").Value = 15.69357920611757
$ws.Range("C4").Value = 8.358064397364053
$ws.Range("E4").Value = 14.21174150322037
$ws.Range("F4").Value = 44.25261380361287
$ws.Range("G4").Value = 48.69725560778203
$ws.Range("H4").Value = 19.38423926403386
$ws.Range("J4").Value = 9.415982752121087
$ws.Range("K4").Value = 11.62533432285984
$ws.Range("L4").Value = 11.62629440247583
$ws.Range("M4").Value = 16.76593501181218
$ws.Range("N4").Value = 22.84436939749173
$ws.Range("B5").Value = 15.65841184789029
$ws.Range("C5").Value = 8.349826534837717
$ws.Range("E5").Value = 14.21628448326076
$ws.Range("F5").Value = 44.26275659832651
$ws.Range("G5").Value = 48.71701380451025
$ws.Range("H5").Value = 19.3973662533933
$ws.Range("J5").Value = 9.417002286430867
$ws.Range("K5").Value = 11.60068344252136
$ws.Range("L5").Value = 11.62814315930414
$ws.Range("M5").Value = 16.76340765192921
$ws.Range("N5").Value = 22.86133766276007
$ws.Range("B6").Value = 15.65261940052942
$ws.Range("C6").Value = 8.3484512243681
$ws.Range("E6").Value = 14.21706193599248
$ws.Range("F6").Value = 44.26457321526468
$ws.Range("G6").Value = 48.7204613952821
$ws.Range("H6").Value = 19.39958812618726
$ws.Range("J6").Value = 9.41717358296542
$ws.Range("K6").Value = 11.59662418336499
$ws.Range("L6").Value = 11.628471742157
$ws.Range("M6").Value = 16.76302722973331
$ws.Range("N6").Value = 22.86418398516883
$ws.Range("B7").Value = 15.69310179475469
$ws.Range("C7").Value = 8.357953796600066
$ws.Range("E7").Value = 14.21180122366939
$ws.Range("F7").Value = 44.25274171597516
$ws.Range("G7").Value = 48.6975108909964
$ws.Range("H7").Value = 19.38441347364265
$ws.Range("J7").Value = 9.415996367682649
$ws.Range("K7").Value = 11.62499960877324
$ws.Range("L7").Value = 11.62631788810685
$ws.Range("M7").Value = 16.76589829805916
$ws.Range("N7").Value = 22.84459631054526
$ws.Range("B8").Value = 15.87832709318426
$ws.Range("C8").Value = 8.398802233916827
$ws.Range("E8").Value = 14.19193875326102
$ws.Range("F8").Value = 44.21977770988585
$ws.Range("G8").Value = 48.62051078994874
$ws.Range("H8").Value = 19.32357895437734
$ws.Range("J8").Value = 9.411091929838911
$ws.Range("K8").Value = 11.75498302037837
$ws.Range("L8").Value = 11.61994646400655
$ws.Range("M8").Value = 16.78424969174372
$ws.Range("N8").Value = 22.76233168469425
$ws.Range("B9").Value = 16.26268236086304
$ws.Range("C9").Value = 8.474954982837053
$ws.Range("E9").Value = 14.16498656190574
$ws.Range("F9").Value = 44.22423587378297
$ws.Range("G9").Value = 48.55703896948171
$ws.Range("H9").Value = 19.22632256561009
$ws.Range("J9").Value = 9.40250937094468
$ws.Range("K9").Value = 12.02529673639212
$ws.Range("L9").Value = 11.61864790803392
$ws.Range("M9").Value = 16.84042362964857
$ws.Range("N9").Value = 22.61592358544795
$ws.Range("B10").Value = 16.55557922147655
$ws.Range("C10").Value = 8.528451676146089
$ws.Range("E10").Value = 14.15245367777285
$ws.Range("F10").Value = 44.26950887390417
$ws.Range("G10").Value = 48.56393620008741
$ws.Range("H10").Value = 19.16831755252309
$ws.Range("J10").Value = 9.396826743915192
$ws.Range("K10").Value = 12.231682333644
$ws.Range("L10").Value = 11.62445681195913
$ws.Range("M10").Value = 16.89376696692195
$ws.Range("N10").Value = 22.51737072336066
$ws.Range("B11").Value = 16.69059012315815
$ws.Range("C11").Value = 8.552246819643768
$ws.Range("E11").Value = 14.14832007060349
$ws.Range("F11").Value = 44.29917714554275
$ws.Range("G11").Value = 48.57871027101633
$ws.Range("H11").Value = 19.14485228824512
$ws.Range("J11").Value = 9.394375138322902
$ws.Range("K11").Value = 12.32690952715592
$ws.Range("L11").Value = 11.62855300708967
$ws.Range("M11").Value = 16.92059881133781
$ws.Range("N11").Value = 22.47447705398957
$ws.Range("B12").Value = 16.74192603970557
$ws.Range("C12").Value = 8.561178787286556
$ws.Range("E12").Value = 14.14697924350472
$ws.Range("F12").Value = 44.31171069514051
$ws.Range("G12").Value = 48.58597693010481
$ws.Range("H12").Value = 19.13638686508169
$ws.Range("J12").Value = 9.393465836304824
$ws.Range("K12").Value = 12.36313243006933
$ws.Range("L12").Value = 11.63031173269107
$ws.Range("M12").Value = 16.93112287315685
$ws.Range("N12").Value = 22.45851189406622
$ws.Range("B13").Value = 16.73086128774706
$ws.Range("C13").Value = 8.559258645600394
$ws.Range("E13").Value = 14.14725804638998
$ws.Range("F13").Value = 44.30895371077767
$ws.Range("G13").Value = 48.58433759461395
$ws.Range("H13").Value = 19.13819134600586
$ws.Range("J13").Value = 9.393660824656367
$ws.Range("K13").Value = 12.35532444047581
$ws.Range("L13").Value = 11.62992375137128
$ws.Range("M13").Value = 16.9288402561333
$ws.Range("N13").Value = 22.46193794027352
$ws.Range("B14").Value = 16.69480960074965
$ws.Range("C14").Value = 8.552983241581311
$ws.Range("E14").Value = 14.14820526662907
$ws.Range("F14").Value = 44.30018229548278
$ws.Range("G14").Value = 48.57927460514554
$ws.Range("H14").Value = 19.14414740676033
$ws.Range("J14").Value = 9.394299948049509
$ws.Range("K14").Value = 12.32988652624771
$ws.Range("L14").Value = 11.62869354683453
$ws.Range("M14").Value = 16.92145738145527
$ws.Range("N14").Value = 22.47315803007897
$ws.Range("B15").Value = 16.67275296207101
$ws.Range("C15").Value = 8.549129092645556
$ws.Range("E15").Value = 14.14881467139408
$ws.Range("F15").Value = 44.29497850179546
$ws.Range("G15").Value = 48.57639106701509
$ws.Range("H15").Value = 19.14785041654865
$ws.Range("J15").Value = 9.39469390915734
$ws.Range("K15").Value = 12.3143253186917
$ws.Range("L15").Value = 11.62796700007381
$ws.Range("M15").Value = 16.91698231670127
$ws.Range("N15").Value = 22.48006679935344
$ws.Range("B16").Value = 16.54678777476495
$ws.Range("C16").Value = 8.526885642561382
$ws.Range("E16").Value = 14.15275528417024
$ws.Range("F16").Value = 44.26775207914949
$ws.Range("G16").Value = 48.56320475609955
$ws.Range("H16").Value = 19.16990987346628
$ws.Range("J16").Value = 9.396989637009241
$ws.Range("K16").Value = 12.22548340152055
$ws.Range("L16").Value = 11.62421824580649
$ws.Range("M16").Value = 16.8920646069374
$ws.Range("N16").Value = 22.52021284640872
$ws.Range("B17").Value = 16.46993411811573
$ws.Range("C17").Value = 8.513100808720054
$ws.Range("E17").Value = 14.15557357927648
$ws.Range("F17").Value = 44.25336956908266
$ws.Range("G17").Value = 48.5580952656485
$ws.Range("H17").Value = 19.18419109612753
$ws.Range("J17").Value = 9.398432084375894
$ws.Range("K17").Value = 12.17130366423292
$ws.Range("L17").Value = 11.62228982694005
$ws.Range("M17").Value = 16.87743171419847
$ws.Range("N17").Value = 22.5453369207374
$ws.Range("B18").Value = 16.42589925098094
$ws.Range("C18").Value = 8.505121295380965
$ws.Range("E18").Value = 14.15734220377301
$ws.Range("F18").Value = 44.24595157601573
$ws.Range("G18").Value = 48.55625186008691
$ws.Range("H18").Value = 19.19268025075784
$ws.Range("J18").Value = 9.399274310144591
$ws.Range("K18").Value = 12.14026899668027
$ws.Range("L18").Value = 11.62131757289042
$ws.Range("M18").Value = 16.8692571081252
$ws.Range("N18").Value = 22.55997014632531
$ws.Range("B19").Value = 16.41102023666649
$ws.Range("C19").Value = 8.502410866873293
$ws.Range("E19").Value = 14.15796641089076
$ws.Range("F19").Value = 44.24358688660914
$ws.Range("G19").Value = 48.55581587507434
$ws.Range("H19").Value = 19.19560175123357
$ws.Range("J19").Value = 9.399561635651931
$ws.Range("K19").Value = 12.12978411730199
$ws.Range("L19").Value = 11.62101194359121
$ws.Range("M19").Value = 16.86653102707566
$ws.Range("N19").Value = 22.56495608470535
$ws.Range("B20").Value = 16.47809811856094
$ws.Range("C20").Value = 8.514573496789756
$ws.Range("E20").Value = 14.15525829600778
$ws.Range("F20").Value = 44.25481221943973
$ws.Range("G20").Value = 48.55852581098081
$ws.Range("H20").Value = 19.18264237493664
$ws.Range("J20").Value = 9.39827723351218
$ws.Range("K20").Value = 12.17705815401589
$ws.Range("L20").Value = 11.62248094981106
$ws.Range("M20").Value = 16.87896441929094
$ws.Range("N20").Value = 22.54264353711963
$ws.Range("B21").Value = 16.70539350316381
$ws.Range("C21").Value = 8.554828624520978
$ws.Range("E21").Value = 14.14792096054338
$ws.Range("F21").Value = 44.30272347544162
$ws.Range("G21").Value = 48.58071636573998
$ws.Range("H21").Value = 19.14238655696159
$ws.Range("J21").Value = 9.394111705491879
$ws.Range("K21").Value = 12.33735409263943
$ws.Range("L21").Value = 11.62904926595899
$ws.Range("M21").Value = 16.92361609309622
$ws.Range("N21").Value = 22.46985488790882
$ws.Range("B22").Value = 16.85514902260531
$ws.Range("C22").Value = 8.580678287563556
$ws.Range("E22").Value = 14.14443355721742
$ws.Range("F22").Value = 44.34160329488854
$ws.Range("G22").Value = 48.60496359872558
$ws.Range("H22").Value = 19.11852734247334
$ws.Range("J22").Value = 9.391500377604315
$ws.Range("K22").Value = 12.44304943765273
$ws.Range("L22").Value = 11.63455126825711
$ws.Range("M22").Value = 16.95491406835347
$ws.Range("N22").Value = 22.42390172892
$ws.Range("B23").Value = 16.77512606837487
$ws.Range("C23").Value = 8.566924192010241
$ws.Range("E23").Value = 14.14617549602288
$ws.Range("F23").Value = 44.32016223808674
$ws.Range("G23").Value = 48.59113153230955
$ws.Range("H23").Value = 19.13103717832157
$ws.Range("J23").Value = 9.392883968347999
$ws.Range("K23").Value = 12.38656255554495
$ws.Range("L23").Value = 11.63150460550059
$ws.Range("M23").Value = 16.93801808763107
$ws.Range("N23").Value = 22.44828003975345
$ws.Range("B24").Value = 16.47440670543576
$ws.Range("C24").Value = 8.513907863999552
$ws.Range("E24").Value = 14.15540037342318
$ws.Range("F24").Value = 44.25415734686361
$ws.Range("G24").Value = 48.5583277532068
$ws.Range("H24").Value = 19.18334168367403
$ws.Range("J24").Value = 9.398347201251793
$ws.Range("K24").Value = 12.17445619180105
$ws.Range("L24").Value = 11.62239411812442
$ws.Range("M24").Value = 16.87827074143202
$ws.Range("N24").Value = 22.54386062683027
$ws.Range("B25").Value = 16.15668362474057
$ws.Range("C25").Value = 8.454781159269558
$ws.Range("E25").Value = 14.17099789177411
$ws.Range("F25").Value = 44.21564146178108
$ws.Range("G25").Value = 48.56481713114655
$ws.Range("H25").Value = 19.2502719963341
$ws.Range("J25").Value = 9.404721202618282
$ws.Range("K25").Value = 11.95068336207681
$ws.Range("L25").Value = 11.61780670047501
$ws.Range("M25").Value = 16.82308873411582
$ws.Range("N25").Value = 22.6539424719451
